$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 245
$ws.Range("E2").Value = 251
$ws.Range("F2").Value = 319
$ws.Range("G2").Value = 208.5500030517578
$ws.Range("H2").Value = 338800000
$ws.Range("I2").Value = "SNOW"

$ws.Range("D3").Value = 245
$ws.Range("E3").Value = 251
$ws.Range("F3").Value = 319
$ws.Range("G3").Value = 208.5500030517578
$ws.Range("H3").Value = 338800000
$ws.Range("I3").Value = "SNOW"

$ws.Range("D4").Value = 245
$ws.Range("E4").Value = 251
$ws.Range("F4").Value = 319
$ws.Range("G4").Value = 208.5500030517578
$ws.Range("H4").Value = 338800000
$ws.Range("I4").Value = "SNOW"

$ws.Range("D5").Value = 245
$ws.Range("E5").Value = 251
$ws.Range("F5").Value = 319
$ws.Range("G5").Value = 208.5500030517578
$ws.Range("H5").Value = 338800000
$ws.Range("I5").Value = "SNOW"

$ws.Range("D6").Value = 245
$ws.Range("E6").Value = 251
$ws.Range("F6").Value = 319
$ws.Range("G6").Value = 208.5500030517578
$ws.Range("H6").Value = 338800000
$ws.Range("I6").Value = "SNOW"

$ws.Range("D7").Value = 254.1000061035156
$ws.Range("E7").Value = 325.8399963378906
$ws.Range("F7").Value = 342
$ws.Range("G7").Value = 229.979995727539
$ws.Range("H7").Value = 338800000
$ws.Range("I7").Value = "SNOW"

$ws.Range("D8").Value = 275.1000061035156
$ws.Range("E8").Value = 259.5400085449219
$ws.Range("F8").Value = 327.4100036621094
$ws.Range("G8").Value = 242.979995727539
$ws.Range("H8").Value = 338800000
$ws.Range("I8").Value = "SNOW"

$ws.Range("D9").Value = 232.3000030517578
$ws.Range("E9").Value = 238.0299987792969
$ws.Range("F9").Value = 246.4750061035156
$ws.Range("G9").Value = 184.7100067138672
$ws.Range("H9").Value = 338800000
$ws.Range("I9").Value = "SNOW"

$ws.Range("D10").Value = 266.0400085449219
$ws.Range("E10").Value = 304.3500061035156
$ws.Range("F10").Value = 307.4200134277344
$ws.Range("G10").Value = 247.8800048828125
$ws.Range("H10").Value = 338800000
$ws.Range("I10").Value = "SNOW"

$ws.Range("D11").Value = 353.8399963378906
$ws.Range("E11").Value = 340.1499938964844
$ws.Range("F11").Value = 405
$ws.Range("G11").Value = 335.0499877929688
$ws.Range("H11").Value = 338800000
$ws.Range("I11").Value = "SNOW"

$ws.Range("D12").Value = 281.7999877929688
$ws.Range("E12").Value = 265.6600036621094
$ws.Range("F12").Value = 329.489990234375
$ws.Range("G12").Value = 240.1000061035156
$ws.Range("H12").Value = 338800000
$ws.Range("I12").Value = "SNOW"

$ws.Range("D13").Value = 170.3099975585938
$ws.Range("E13").Value = 127.6500015258789
$ws.Range("F13").Value = 187.229995727539
$ws.Range("G13").Value = 112.0999984741211
$ws.Range("H13").Value = 338800000
$ws.Range("I13").Value = "SNOW"

$ws.Range("D14").Value = 146.6499938964844
$ws.Range("E14").Value = 180.9499969482422
$ws.Range("F14").Value = 205.6600036621093
$ws.Range("G14").Value = 143.0399932861328
$ws.Range("H14").Value = 338800000
$ws.Range("I14").Value = "SNOW"

$ws.Range("D15").Value = 164.7400054931641
$ws.Range("E15").Value = 142.8999938964844
$ws.Range("F15").Value = 168
$ws.Range("G15").Value = 122.7699966430664
$ws.Range("H15").Value = 338800000
$ws.Range("I15").Value = "SNOW"

$ws.Range("D16").Value = 158.5899963378906
$ws.Range("E16").Value = 154.3800048828125
$ws.Range("F16").Value = 178.6999969482422
$ws.Range("G16").Value = 145.2799987792969
$ws.Range("H16").Value = 338800000
$ws.Range("I16").Value = "SNOW"

$ws.Range("D17").Value = 148.1499938964844
$ws.Range("E17").Value = 165.3600006103516
$ws.Range("F17").Value = 185
$ws.Range("G17").Value = 142.4420013427734
$ws.Range("H17").Value = 338800000
$ws.Range("I17").Value = "SNOW"

$ws.Range("D18").Value = 175.2299957275391
$ws.Range("E18").Value = 156.8500061035156
$ws.Range("F18").Value = 177.7100067138672
$ws.Range("G18").Value = 142.3999938964844
$ws.Range("H18").Value = 338800000
$ws.Range("I18").Value = "SNOW"

$ws.Range("D19").Value = 145.3000030517578
$ws.Range("E19").Value = 187.6799926757812
$ws.Range("F19").Value = 192.6600036621093
$ws.Range("G19").Value = 140.1399993896484
$ws.Range("H19").Value = 338800000
$ws.Range("I19").Value = "SNOW"

$ws.Range("D20").Value = 197
$ws.Range("E20").Value = 188.2799987792969
$ws.Range("F20").Value = 237.7200012207031
$ws.Range("G20").Value = 180.6799926757812
$ws.Range("H20").Value = 338800000
$ws.Range("I20").Value = "SNOW"

$ws.Range("D21").Value = 158.1399993896484
$ws.Range("E21").Value = 136.1799926757812
$ws.Range("F21").Value = 168.8000030517578
$ws.Range("G21").Value = 133.5899963378906
$ws.Range("H21").Value = 338800000
$ws.Range("I21").Value = "SNOW"

$ws.Range("D22").Value = 129.3999938964844
$ws.Range("E22").Value = 114.2300033569336
$ws.Range("F22").Value = 135.7200012207031
$ws.Range("G22").Value = 107.9300003051758
$ws.Range("H22").Value = 338800000
$ws.Range("I22").Value = "SNOW"

$ws.Range("D23").Value = 115.9800033569336
$ws.Range("E23").Value = 174.8000030517578
$ws.Range("F23").Value = 178.6999969482422
$ws.Range("G23").Value = 113.2300033569336
$ws.Range("H23").Value = 338800000
$ws.Range("I23").Value = "SNOW"

$ws.Range("D24").Value = 176.1699981689453
$ws.Range("E24").Value = 177.1000061035156
$ws.Range("F24").Value = 194.3999938964844
$ws.Range("G24").Value = 160.6600036621094
$ws.Range("H24").Value = 338800000
$ws.Range("I24").Value = "SNOW"

$ws.Range("D25").Value = 168.3950042724609
$ws.Range("E25").Value = 205.6699981689453
$ws.Range("F25").Value = 209.8000030517578
$ws.Range("G25").Value = 164.2400054931641
$ws.Range("H25").Value = 338800000
$ws.Range("I25").Value = "SNOW"

$ws.Range("D26").Value = 216.8000030517578
$ws.Range("E26").Value = 238.6600036621093
$ws.Range("F26").Value = 249.9900054931641
$ws.Range("G26").Value = 188
$ws.Range("H26").Value = 338800000
$ws.Range("I26").Value = "SNOW"
